# Completed Global filter and keyword search filter functionalities for All
# Locations tab in Manage Pages page. Done Exporting data and comparing data
# between UI and excel.
#
# Adds a new "Configuration" worksheet (after the existing "TPSEE" sheet)
# that holds the Abu Dhabi / AE "Country_A" filter-scenario row used by the
# new All-Locations keyword search, and updates the sheet selections so the
# new sheet becomes the active tab.

$wb = $excel.ActiveWorkbook
$tpsee = $wb.Worksheets.Item(1)

# --- Add the new "Configuration" sheet right after "TPSEE" ---------------
$cfg = $wb.Worksheets.Add($null, $tpsee)
$cfg.Name = "Configuration"

# --- Header row (same headers as TPSEE) -----------------------------------
$cfg.Range("A1").Value = "FilterScenarioNum"
$cfg.Range("B1").Value = "Group"
$cfg.Range("C1").Value = "Country"
$cfg.Range("D1").Value = "State"
$cfg.Range("E1").Value = "City"
$cfg.Range("F1").Value = "Location"

# --- Data row: new Abu Dhabi / UAE test location --------------------------
$cfg.Range("A2").Value = "Scenario1"
$cfg.Range("B2").Value = "Country_A"
$cfg.Range("C2").Value = "AE"
$cfg.Range("D2").Value = "AUH"
$cfg.Range("E2").Value = "Abu Dhabi"
$cfg.Range("F2").Value = "TestLocation5, Theyab Bin Eissa St, +971600522252"

# Match the Consolas "code" look used for Country/State/City on TPSEE by
# copying the formatting from that sheet onto the same cells here.
$tpsee.Range("C2:E2").Copy()
$cfg.Range("C2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths (best-fit look, approximated in character units) -------
$cfg.Columns.Item(1).ColumnWidth = 16.916666666666668
$cfg.Columns.Item(2).ColumnWidth = 9.416666666666666
$cfg.Columns.Item(5).ColumnWidth = 9.083333333333334
$cfg.Columns.Item(6).ColumnWidth = 44.916666666666664

# --- Selections -------------------------------------------------------------
# TPSEE is no longer the active tab; its whole data range is selected.
$tpsee.Range("A1:F2").Select()

# Configuration becomes the active / selected tab, with F3 as the active cell.
$cfg.Select()
$cfg.Range("F3").Select()

Write-Output "Added Configuration sheet"
